$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format to prevent Excel from
# auto-converting numeric-looking strings (e.g. "1.00" -> 1) and
# losing formatting such as trailing zeros or thousands separators.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '72.036.25'
$ws.Range("E2").Value = '  +5.97%  '

# Row 3
$ws.Range("D3").Value = '2.531.13'
$ws.Range("E3").Value = '  +4.99%  '

# Row 4
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").Value = '577.86'
$ws.Range("E5").Value = '  +4.17%  '

# Row 6
$ws.Range("D6").Value = '176.91'
$ws.Range("E6").Value = '  +11.32%  '

# Row 7
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.14%  '

# Row 8
$ws.Range("D8").Value = '0.523'
$ws.Range("E8").Value = '  +3.53%  '

# Row 9
$ws.Range("D9").Value = '0.186'
$ws.Range("E9").Value = '  +14.15%  '

# Row 10
$ws.Range("D10").Value = '2.528.93'
$ws.Range("E10").Value = '  +5.00%  '

# Row 11
$ws.Range("E11").Value = '  -0.93%  '

# Row 12
$ws.Range("D12").Value = '0.345'
$ws.Range("E12").Value = '  +4.88%  '

# Row 13
$ws.Range("D13").Value = '4.75'
$ws.Range("E13").Value = '  +1.55%  '

# Row 14
$ws.Range("D14").Value = '0.0000187'
$ws.Range("E14").Value = '  +7.61%  '

# Row 15
$ws.Range("D15").Value = '71.930.00'
$ws.Range("E15").Value = '  +6.03%  '

# Row 16
$ws.Range("D16").Value = '2.992.58'
$ws.Range("E16").Value = '  +4.79%  '

# Row 17
$ws.Range("D17").Value = '25.12'
$ws.Range("E17").Value = '  +9.98%  '

# Row 18
$ws.Range("D18").Value = '2.521.99'
$ws.Range("E18").Value = '  +4.52%  '

# Row 19
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = '11.32'
$ws.Range("E19").Value = '  +9.58%  '

# Row 20
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '7.74'
$ws.Range("E20").Value = '  +12.88%  '

# Row 21
$ws.Range("D21").Value = '354.18'
$ws.Range("E21").Value = '  +7.38%  '

# Row 22
$ws.Range("D22").Value = '2.17'
$ws.Range("E22").Value = '  +16.53%  '

# Row 23
$ws.Range("D23").Value = '3.97'
$ws.Range("E23").Value = '  +4.83%  '

# Row 24
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  +0.09%  '

# Row 25
$ws.Range("D25").Value = '67.65'
$ws.Range("E25").Value = '  +2.52%  '

# Row 26
$ws.Range("D26").Value = '3.97'
$ws.Range("E26").Value = '  +8.99%  '

# Row 27
$ws.Range("D27").Value = '8.87'
$ws.Range("E27").Value = '  +8.99%  '

# Row 28
$ws.Range("D28").Value = '2.634.83'
$ws.Range("E28").Value = '  +3.99%  '

# Row 29
$ws.Range("D29").Value = '0.995'
$ws.Range("E29").Value = '  -1.04%  '

# Row 30
$ws.Range("D30").Value = '0.0₃0907'
$ws.Range("E30").Value = '  +12.30%  '

# Row 31
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").Value = '493.42'
$ws.Range("E31").Value = '  +17.76%  '

# Row 32
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '7.63'
$ws.Range("E32").Value = '  +7.89%  '

# Row 33
$ws.Range("D33").Value = '1.31'
$ws.Range("E33").Value = '  +15.39%  '

# Row 34
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  +0.00%  '

# Row 35
$ws.Range("D35").Value = '1.68'
$ws.Range("E35").Value = '  +4.82%  '

# Row 36
$ws.Range("D36").Value = '158.99'
$ws.Range("E36").Value = '  -1.14%  '

# Row 37
$ws.Range("D37").Value = '0.116'
$ws.Range("E37").Value = '  +10.01%  '

# Row 38
$ws.Range("D38").Value = '19.28'
$ws.Range("E38").Value = '  +1.52%  '

# Row 39
$ws.Range("D39").Value = '18.74'
$ws.Range("E39").Value = '  +5.41%  '

# Row 40
$ws.Range("E40").Value = '  +0.04%  '

# Row 41
$ws.Range("D41").Value = '4.69'
$ws.Range("E41").Value = '  +9.09%  '

# Row 42
$ws.Range("D42").Value = '1.60'
$ws.Range("E42").Value = '  +8.59%  '

# Row 43
$ws.Range("D43").Value = '0.312'
$ws.Range("E43").Value = '  +5.57%  '

# Row 44
$ws.Range("D44").Value = '38.23'
$ws.Range("E44").Value = '  +2.35%  '

# Row 45
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '2.27'
$ws.Range("E45").Value = '  +15.53%  '

# Row 46
$ws.Range("B46").Value = 'ImmutableX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D46").Value = '1.13'
$ws.Range("E46").Value = '  +5.85%  '

# Row 47
$ws.Range("D47").Value = '144.94'
$ws.Range("E47").Value = '  +10.71%  '

# Row 48
$ws.Range("D48").Value = '3.51'
$ws.Range("E48").Value = '  +6.41%  '

# Row 49
$ws.Range("D49").Value = '0.514'
$ws.Range("E49").Value = '  +7.53%  '

# Row 50
$ws.Range("D50").Value = '0.0740'
$ws.Range("E50").Value = '  +4.27%  '

# Row 51
$ws.Range("D51").Value = '0.577'
$ws.Range("E51").Value = '  +4.44%  '
